# linkdata.xlsx edit: update Mars 6U CubeSat (column E) link-budget inputs,
# add threaded-comment discussion about the BIRD example column data, and
# restore the view/selection state left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Mars 6U CubeSat (column E) parameter updates -------------------------
$ws.Range("E3").Value = 8       # S/C transmitter power
$ws.Range("E4").Value = 1000    # Ground Station transmitter power
$ws.Range("E6").Value = 0.8     # L receiver
$ws.Range("E10").Value = 60     # Antenna D ground station
$ws.Range("E17").Value = 0.5    # Payload pixel size
$ws.Range("E20").Value = 18     # Payload downlink time

# --- Threaded comments ------------------------------------------------------
# Reply to the existing F16 thread (OG:45 / BIRD: 1) explaining the BIRD
# payload mock data.
$f16 = $ws.Range("F16").CommentThreaded
$f16.AddReply("BIRD payload data is mock data to generate 1e6 bit/s, so the actual dimensions are bullshit, they only amount to 1e6 in total")

# New thread on H16 (BIRD example column, Payload swath width angle) noting
# that the data is mock data, with a reply clarifying why.
$h16 = $ws.Range("H16").AddCommentThreaded("This data is bullshit")
$h16.AddReply("It was made to give a data rate of 1e6 bit/s")

# --- View / selection state --------------------------------------------------
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H16").Select()
$ws.Range("E4").Select()
